# edit.ps1
# Applies the "dictionnaire des donnees en cours (world idee)" change:
#  1. Removes the "Lancement..." / Equipe_id / Equipe_nom block from just
#     below the "Idee" heading (collapsing those 5 paragraphs into one
#     empty paragraph).
#  2. Moves <w:lastRenderedPageBreak/> from the "Questions" heading run to
#     the "Supprimer une question" bullet run.
#  3. Re-inserts the "Lancement..." / Equipe_id / Equipe_nom block (plus a
#     large new set of data-dictionary paragraphs: theme_id, theme_nom,
#     Question_id, Question_niveau, Reponse_id, Reponse_resutla) at the
#     very end of the document, right before the final section break.

$d = $word.ActiveDocument

function Get-ParaIndexByText($text) {
    $cnt = $d.Paragraphs.Count
    for ($i = 1; $i -le $cnt; $i++) {
        if ($d.Paragraphs.Item($i).Range.Text.Contains($text)) {
            return $i
        }
    }
    return -1
}

$pkgOpen = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">'
$pkgClose = '</w:document></pkg:xmlData></pkg:part></pkg:package>'

# ---------------------------------------------------------------------
# Step 1: collapse the "Lancement.../Equipe_id/Equipe_nom" block (right
# after the "Idee" heading, including the blank paragraph that trails it)
# down to a single empty paragraph.
# ---------------------------------------------------------------------
$startIdx = Get-ParaIndexByText("Lancement d")
$endIdx = Get-ParaIndexByText("Equipe_nom")

$startPos = $d.Paragraphs.Item($startIdx).Range.Start
$endPos = $d.Paragraphs.Item($endIdx + 1).Range.End

$blockRange = $d.Range($startPos, $endPos)
$emptyParaXml = $pkgOpen + '<w:body><w:p><w:pPr><w:spacing w:after="0"/></w:pPr></w:p></w:body>' + $pkgClose
$blockRange.InsertXML($emptyParaXml)

# ---------------------------------------------------------------------
# Step 2a: remove <w:lastRenderedPageBreak/> from the "Questions" heading.
# ---------------------------------------------------------------------
$qIdx = Get-ParaIndexByText("Questions")
$qPara = $d.Paragraphs.Item($qIdx)
$qXml = $pkgOpen + '<w:body><w:p><w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Questions</w:t></w:r></w:p></w:body>' + $pkgClose
$qPara.Range.InsertXML($qXml)

# ---------------------------------------------------------------------
# Step 2b: add <w:lastRenderedPageBreak/> to the "Supprimer une question"
# bullet item.
# ---------------------------------------------------------------------
$sIdx = Get-ParaIndexByText("Supprimer une question")
$sPara = $d.Paragraphs.Item($sIdx)
$sXml = $pkgOpen + '<w:body><w:p><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="4"/></w:numPr></w:pPr><w:r><w:lastRenderedPageBreak/><w:t>Supprimer une question</w:t></w:r></w:p></w:body>' + $pkgClose
$sPara.Range.InsertXML($sXml)

# ---------------------------------------------------------------------
# Step 3: append the data-dictionary paragraphs at the very end of the
# document body, right before the final section break.
# ---------------------------------------------------------------------
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$newBlock = '<w:p><w:pPr><w:spacing w:after="0"/></w:pPr><w:r><w:t>Lancement d&#8217;une partie de quiz param&#233;tr&#233;e (choix du quiz, nom des &#233;quipes)</w:t></w:r></w:p>' +
  '<w:p><w:pPr><w:spacing w:after="0"/></w:pPr></w:p>' +
  '<w:p><w:pPr><w:spacing w:after="0"/></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>Equipe_id</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>' +
  '<w:p><w:pPr><w:spacing w:after="0"/></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>Equipe_nom</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>' +
  '<w:p><w:pPr><w:spacing w:after="0"/></w:pPr></w:p>' +
  '<w:p><w:pPr><w:spacing w:after="0"/></w:pPr><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>t</w:t></w:r><w:r><w:t>heme</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>_id</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>' +
  '<w:p><w:pPr><w:spacing w:after="0"/></w:pPr><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>t</w:t></w:r><w:r><w:t>heme</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>_nom</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>' +
  '<w:p><w:pPr><w:spacing w:after="0"/></w:pPr></w:p>' +
  '<w:p><w:pPr><w:spacing w:after="0"/></w:pPr></w:p>' +
  '<w:p><w:pPr><w:spacing w:after="0"/></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>Question_id</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>' +
  '<w:p><w:pPr><w:spacing w:after="0"/></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>Question_</w:t></w:r><w:r><w:t>niveau</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r></w:p>' +
  '<w:p><w:pPr><w:spacing w:after="0"/></w:pPr></w:p>' +
  '<w:p><w:pPr><w:spacing w:after="0"/></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>Reponse_id</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>' +
  '<w:p><w:pPr><w:spacing w:after="0"/></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>Reponse_res</w:t></w:r><w:r><w:t>utla</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r></w:p>' +
  '<w:p><w:pPr><w:spacing w:after="0"/></w:pPr></w:p>' +
  '<w:p><w:pPr><w:spacing w:after="0"/></w:pPr></w:p>'

$finalXml = $pkgOpen + '<w:body>' + $newBlock + '</w:body>' + $pkgClose
$lastPara.Range.InsertXML($finalXml)

Write-Output "edit applied"
